$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $src = $ws.Range("A1")
    $src.NumberFormat = "@"
    $src.Value2 = $val
    $src.Copy()
    $dst = $ws.Range($cellRef)
    $dst.PasteSpecial(-4163, 0, $false, $false)
    $src.Clear()
}

Set-TextValue "D2" '76.503.57'
Set-TextValue "E2" '  +0.68%  '

Set-TextValue "D3" '2.958.43'
Set-TextValue "E3" '  +2.16%  '

Set-TextValue "E4" '  -0.03%  '

Set-TextValue "D5" '199.90'
Set-TextValue "E5" '  +1.82%  '

Set-TextValue "D6" '597.66'
Set-TextValue "E6" '  -0.15%  '

Set-TextValue "E7" '  -0.04%  '

Set-TextValue "D8" '0.554'
Set-TextValue "E8" '  -0.14%  '

Set-TextValue "D9" '0.202'
Set-TextValue "E9" '  +4.57%  '

Set-TextValue "D10" '2.956.72'
Set-TextValue "E10" '  +2.06%  '

Set-TextValue "D11" '0.444'
Set-TextValue "E11" '  +9.88%  '

Set-TextValue "E12" '  +0.39%  '

Set-TextValue "E13" '  +0.36%  '

Set-TextValue "D14" '3.501.68'
Set-TextValue "E14" '  +2.19%  '

Set-TextValue "D15" '28.67'
Set-TextValue "E15" '  +4.70%  '

Set-TextValue "D16" '76.324.69'
Set-TextValue "E16" '  +0.58%  '

Set-TextValue "E17" '  -0.88%  '

Set-TextValue "D18" '2.954.59'
Set-TextValue "E18" '  +2.41%  '

Set-TextValue "D19" '13.73'
Set-TextValue "E19" '  +8.82%  '

Set-TextValue "D20" '8.84'
Set-TextValue "E20" '  -0.90%  '

Set-TextValue "D21" '377.41'
Set-TextValue "E21" '  -0.07%  '

Set-TextValue "E22" '  -1.07%  '

Set-TextValue "E23" '  +3.97%  '

Set-TextValue "D24" '72.80'
Set-TextValue "E24" '  +1.81%  '

Set-TextValue "D25" '0.999'
Set-TextValue "E25" '  -0.40%  '

Set-TextValue "E26" '  +2.42%  '

Set-TextValue "D27" '4.34'
Set-TextValue "E27" '  +2.53%  '

Set-TextValue "D28" '9.72'
Set-TextValue "E28" '  -1.04%  '

Set-TextValue "E29" '  -1.13%  '

Set-TextValue "D30" '0.996'
Set-TextValue "E30" '  -0.35%  '

Set-TextValue "D31" '8.58'
Set-TextValue "E31" '  +10.19%  '

Set-TextValue "D32" '1.40'
Set-TextValue "E32" '  -0.89%  '

Set-TextValue "D33" '498.32'
Set-TextValue "E33" '  -1.73%  '

Set-TextValue "E34" '  +1.06%  '

Set-TextValue "D35" '0.999'
Set-TextValue "E35" '  +0.02%  '

Set-TextValue "D36" '166.02'
Set-TextValue "E36" '  +0.62%  '

Set-TextValue "D37" '20.41'
Set-TextValue "E37" '  +0.68%  '

Set-TextValue "D38" '0.392'
Set-TextValue "E38" '  +13.48%  '

Set-TextValue "E39" '  +18.14%  '

Set-TextValue "D40" '19.98'
Set-TextValue "E40" '  +1.42%  '

Set-TextValue "E41" '  -2.25%  '

Set-TextValue "B42" 'Aave'
Set-TextValue "C42" 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue "D42" '181.78'
Set-TextValue "E42" '  -0.87%  '

Set-TextValue "B43" 'USDe'
Set-TextValue "C43" 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue "D43" '1.00'
Set-TextValue "E43" '  +0.00%  '

Set-TextValue "D44" '4.97'
Set-TextValue "E44" '  -1.11%  '

Set-TextValue "E45" '  -1.36%  '

Set-TextValue "E46" '  -1.87%  '

Set-TextValue "D47" '39.81'
Set-TextValue "E47" '  -1.46%  '

Set-TextValue "E48" '  +1.81%  '

Set-TextValue "D49" '3.91'
Set-TextValue "E49" '  +3.65%  '

Set-TextValue "E50" '  -2.36%  '

Set-TextValue "D51" '0.670'
Set-TextValue "E51" '  +0.05%  '
